$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "daf"
$ws.Range("B19").Value = "df"
$ws.Range("C19").Value = "fds"
$ws.Range("D19").Value = "df"
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""

$ws.Range("A20").Value = "dsf"
$ws.Range("B20").Value = "dsf"
$ws.Range("C20").Value = "ds"
$ws.Range("D20").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
